$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set header labels
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Fill data for columns I (I0) and J (IF), rows 2-62
$data = @(
    @(7,8),
    @(7,7),
    @(6,6),
    @(8,8),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(9,9),
    @(9,9),
    @(7,8),
    @(7,7),
    @(7,7),
    @(7,8),
    @(9,9),
    @(8,8),
    @(7,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,8),
    @(8,8),
    @(7,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,9),
    @(8,8),
    @(7,7),
    @(9,9),
    @(10,10),
    @(8,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,8),
    @(9,9),
    @(9,9),
    @(9,10),
    @(7,8),
    @(9,9),
    @(7,8),
    @(8,8),
    @(9,9),
    @(7,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(7,8),
    @(9,9),
    @(6,6),
    @(3,3),
    @(8,8),
    @(4,4),
    @(5,5)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}

Write-Host "Applied I0/IF columns."
